$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B, C, D, E are treated as text so Excel does not
# auto-convert numeric-looking / percentage-looking strings into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '57.912.31'
$ws.Range('E2').Value = '  -3.97%  '
$ws.Range('D3').Value = '3.109.73'
$ws.Range('E3').Value = '  -5.75%  '
$ws.Range('D5').Value = '519.90'
$ws.Range('E5').Value = '  -6.71%  '
$ws.Range('D6').Value = '130.91'
$ws.Range('E6').Value = '  -6.99%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.110.34'
$ws.Range('E8').Value = '  -5.76%  '
$ws.Range('D9').Value = '0.441'
$ws.Range('E9').Value = '  -5.53%  '
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  -8.75%  '
$ws.Range('E11').Value = '  -9.99%  '
$ws.Range('D12').Value = '0.369'
$ws.Range('E12').Value = '  -9.32%  '
$ws.Range('D13').Value = '3.651.99'
$ws.Range('E13').Value = '  -5.50%  '
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').Value = '24.76'
$ws.Range('E15').Value = '  -6.76%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '57.954.33'
$ws.Range('E16').Value = '  -3.87%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.116.91'
$ws.Range('E17').Value = '  -5.41%  '
$ws.Range('D18').Value = '0.0000149'
$ws.Range('E18').Value = '  -9.13%  '
$ws.Range('D19').Value = '5.64'
$ws.Range('E19').Value = '  -7.15%  '
$ws.Range('D20').Value = '12.82'
$ws.Range('E20').Value = '  -6.09%  '
$ws.Range('D21').Value = '7.77'
$ws.Range('E21').Value = '  -8.78%  '
$ws.Range('D22').Value = '336.87'
$ws.Range('E22').Value = '  -9.83%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '0.504'
$ws.Range('E24').Value = '  -5.10%  '
$ws.Range('D25').Value = '66.53'
$ws.Range('E25').Value = '  -7.74%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  -3.88%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0911'
$ws.Range('E28').Value = '  -11.00%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '6.69'
$ws.Range('E30').Value = '  -4.68%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.26'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.84'
$ws.Range('E32').Value = '  -8.76%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '6.77'
$ws.Range('E33').Value = '  -7.10%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '21.14'
$ws.Range('E34').Value = '  -6.21%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = '158.28'
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '4.72'
$ws.Range('E36').Value = '  -5.97%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '6.10'
$ws.Range('E37').Value = '  -7.83%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '1.35'
$ws.Range('E38').Value = '  -11.20%  '
$ws.Range('B39').Value = 'RenzoRestakedETH'
$ws.Range('C39').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D39').Value = '3.145.73'
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '40.29'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.0673'
$ws.Range('E41').Value = '  -6.79%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '23.08'
$ws.Range('E42').Value = '  -8.24%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.684'
$ws.Range('E43').Value = '  -8.45%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').Value = '1.06'
$ws.Range('E44').Value = '  -5.49%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '3.88'
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.260.74'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.42'
$ws.Range('E48').Value = '  -9.64%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '6.09'
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '20.18'
$ws.Range('E50').Value = '  -5.85%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0230'
$ws.Range('E51').Value = '  -7.51%  '
